$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13082
$ws1.Range("F5").Value = 92
$ws1.Range("F10").Value = 13042
$ws1.Range("F11").Value = 300
$ws1.Range("F13").Value = 8737
$ws1.Range("F14").Value = 7773

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13082
$ws4.Range("F6").Value = 92
$ws4.Range("F11").Value = 13042
$ws4.Range("F12").Value = 300
$ws4.Range("F14").Value = 8737
$ws4.Range("F15").Value = 7773
